$d = $word.ActiveDocument

# --- 1. Remove the stray "_GoBack" bookmark that sits, on its own, in the
#        empty paragraph right before the "<< Installation Guide >>"
#        heading. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Update the installation instructions:
#        'Execute the "ConnectFour.bat" inside the build folder'
#     -> 'Execute the "PlayGame.cmd" inside the build folder'
#     and leave a "_GoBack" bookmark right after "cmd" (where the cursor
#     was left after the edit), same as the authored change. ---

# Locate the paragraph that holds the instruction text.
$para = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i).Range
    if ($candidate.Text -like '*ConnectFour.bat*') {
        $para = $candidate
        break
    }
}

$r = $para.Duplicate
$r.Find.Execute("ConnectFour", $false, $false, $false, $false, $false, $true, 1, $false, "PlayGame", 2)

$r = $para.Duplicate
$r.Find.Execute("bat", $false, $false, $false, $false, $false, $true, 1, $false, "cmd", 2)

# The two replacements above merge the whole sentence back into one run
# (same rPr throughout). Re-split it into the same runs the source edit
# produced - 'Execute the "' | 'PlayGame' | '.' | 'cmd' - using
# throw-away bookmarks to pin each boundary, then discarding them.
$b1 = $para.Duplicate
$b1.Find.Execute("PlayGame")
$b1.Collapse(1)
$d.Bookmarks.Add("_TmpBoundary1", $b1)

$b2 = $para.Duplicate
$b2.Find.Execute("PlayGame")
$b2.Collapse(0)
$d.Bookmarks.Add("_TmpBoundary2", $b2)

$b3 = $para.Duplicate
$b3.Find.Execute(".")
$b3.Collapse(0)
$d.Bookmarks.Add("_TmpBoundary3", $b3)

$d.Bookmarks.Item("_TmpBoundary1").Delete()
$d.Bookmarks.Item("_TmpBoundary2").Delete()
$d.Bookmarks.Item("_TmpBoundary3").Delete()

# Finally, drop the real "_GoBack" bookmark right after "cmd", splitting
# off the trailing '" inside the build folder' run.
$goBack = $para.Duplicate
$goBack.Find.Execute("cmd")
$goBack.Collapse(0)
$d.Bookmarks.Add("_GoBack", $goBack)
